# Apply updated "dSF" (column F) values after repulling/pushing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -4
    8  = -5
    9  = 11
    13 = -6
    15 = 0
    16 = 2
    20 = -1
    22 = -6
    28 = -2
    29 = -6
    30 = -2
    33 = -7
    36 = -3
    37 = -3
    43 = 3
    44 = 0
    47 = -1
    48 = -4
    50 = -2
    58 = -3
    60 = -4
    61 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
